$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 388; this shifts the existing rows
# 388-422 down to 389-423 (preserving their data), and leaves a fresh
# blank row 388 for the new record.
$ws.Rows.Item(388).Insert()

# Populate the new row 388 with the new weekly record.
$ws.Range("A388").Value = 3
$ws.Range("B388").Value = "Femacal de La Calera"
$ws.Range("C388").Value = "Coquimbo"
$ws.Range("D388").Value = 44578
$ws.Range("E388").Value = 5
$ws.Range("F388").Value = 100112045
$ws.Range("G388").Value = "Zapallo"
$ws.Range("H388").Value = "Camote"
$ws.Range("I388").Value = "1a nueva(o)"
$ws.Range("J388").Value = 510
$ws.Range("K388").Value = 550
$ws.Range("L388").Value = 600
$ws.Range("M388").Value = 575
$ws.Range("N388").Value = "$/kilo (volumen en unidades)"
$ws.Range("O388").Value = "Provincia de Talca"
$ws.Range("P388").Value = 575
$ws.Range("Q388").Value = 1
$ws.Range("R388").Value = "Hortaliza"

# Match the date-format style used by the other rows' Fecha column.
$ws.Range("D388").NumberFormat = $ws.Range("D389").NumberFormat
